# Monitoreo CRGs - Incorporacion de Importe CRG al monitoreo de CRGs facturados.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update connection parameters on the "Tabla1" parameter table:
#  - workdirectory (row 5) now points at the new shared/synced folder location
#  - user (row 3) and password (row 4) now both use the new "odoo" credential value
$ws.Range("B5").Value = "C:/Users/iachenbach/Gobierno de la Ciudad de Buenos Aires/Pablo Alfredo Gadea - Tablero Facoep P BI/FACOEP/DBA/Reportes BI/2021/Monitoreo CRGs/"
$ws.Range("B3").Value = "odoo"
$ws.Range("B4").Value = "odoo"

# Leave the active cell where the author left it while working further down the sheet
$ws.Range("B14").Select()
